# ---------------------------------------------------------------------------
# Lab3A instructions edit:
#   1. After "Using Tag Helpers" add a tab-only run, then a run with a tab
#      followed by "(use at least 2 kinds)".
#   2. After "Using Unit tests" add a tab-only run, then a run with a tab
#      followed by "(write at least 3 tests)", and move the "_GoBack"
#      bookmark to sit at the end of that paragraph.
#   3. Remove the "_GoBack" bookmark that used to sit after the
#      "wins and losses" bullet (it moved, per step 2).
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

function Append-FormattedRun {
    param($TargetParagraphIndex, $SourceSearchText, $NewText, $ScratchParagraphIndex)

    # Locate the run whose character formatting we want to clone.
    $srcParagraph = $d.Paragraphs.Item($TargetParagraphIndex)
    $srcRange = $srcParagraph.Range.Duplicate
    $srcRange.Find.Execute($SourceSearchText)
    $srcLen = $srcRange.End - $srcRange.Start
    $srcRange.Copy()

    # Paste a same-formatted copy into a scratch (empty) paragraph so we can
    # safely rewrite its text without Word merging it into a neighboring run.
    $scratchParagraph = $d.Paragraphs.Item($ScratchParagraphIndex)
    $scratchInsertPoint = $d.Range($scratchParagraph.Range.Start, $scratchParagraph.Range.Start)
    $scratchStart = $scratchInsertPoint.Start
    $scratchInsertPoint.Paste()
    $scratchRange = $d.Range($scratchStart, $scratchStart + $srcLen)
    $scratchRange.Text = $NewText
    $scratchFinal = $d.Range($scratchStart, $scratchStart + $NewText.Length)

    # Copy the rewritten, correctly-formatted run back out ...
    $scratchFinal.Copy()

    # ... and paste it at the end of the target paragraph (just before the
    # paragraph mark) so it becomes a new trailing run.
    $targetParagraph = $d.Paragraphs.Item($TargetParagraphIndex)
    $targetPos = $targetParagraph.Range.End - 1
    $targetRange = $d.Range($targetPos, $targetPos)
    $targetRange.Paste()

    # Clean the scratch paragraph back to empty.
    $scratchParagraph2 = $d.Paragraphs.Item($ScratchParagraphIndex)
    if (($scratchParagraph2.Range.End - 1) -gt $scratchParagraph2.Range.Start) {
        $scratchClean = $d.Range($scratchParagraph2.Range.Start, $scratchParagraph2.Range.End - 1)
        $scratchClean.Delete()
    }
}

$TAB = [char]9

# --- 1. "Using Tag Helpers" bullet ----------------------------------------
Append-FormattedRun 3 "Tag Helpers" "$TAB" 6
Append-FormattedRun 3 "Tag Helpers" "$TAB(use at least 2 kinds)" 6

# --- 2. "Using Unit tests" bullet ------------------------------------------
Append-FormattedRun 5 "Unit tests" "$TAB" 6
Append-FormattedRun 5 "Unit tests" "$TAB(write at least 3 tests)" 6

# --- 3. Move the "_GoBack" bookmark from the "wins and losses" bullet to
#        the end of the "Using Unit tests" paragraph. ----------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$unitTestsParagraph = $d.Paragraphs.Item(5)
$bookmarkPos = $unitTestsParagraph.Range.End - 1
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
